$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/URL/percentage updates (unambiguous text, Excel keeps them as text).
$ws.Range('D2').Value = '25.965.78'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').Value = '1.757.07'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  +0.38%  '
$ws.Range('E5').Value = '  -1.11%  '
$ws.Range('E6').Value = '  +0.50%  '
$ws.Range('E7').Value = '  +3.42%  '
$ws.Range('B8').Value = 'OKB'
$ws.Range('C8').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('E8').Value = '  -3.18%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('E9').Value = '  +2.23%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('E10').Value = '  +0.45%  '
$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').Value = '1.758.14'
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('E12').Value = '  +1.02%  '
$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('E13').Value = '  -0.07%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('E14').Value = '  +7.77%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('E15').Value = '  +0.02%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('E16').Value = '  +1.15%  '
$ws.Range('B17').Value = 'BinanceUSD'
$ws.Range('C17').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('E17').Value = '  +0.32%  '
$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('E18').Value = '  +0.23%  '
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '25.971.50'
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('E20').Value = '  +0.07%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('E21').Value = '  -0.62%  '
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '1.978.94'
$ws.Range('E22').Value = '  +0.56%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('E24').Value = '  +3.69%  '
$ws.Range('B25').Value = 'Chainlink'
$ws.Range('C25').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('E25').Value = '  +1.65%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('E26').Value = '  -0.74%  '
$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('E27').Value = '  -3.32%  '
$ws.Range('E28').Value = '  -0.55%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('E29').Value = '  +1.81%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('E30').Value = '  +0.76%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('E31').Value = '  +3.62%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('E32').Value = '  -1.64%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('E33').Value = '  -1.65%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('E34').Value = '  -1.73%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('E35').Value = '  +1.06%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('E36').Value = '  +0.74%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('E37').Value = '  +0.42%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('E38').Value = '  +1.53%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('E39').Value = '  +2.88%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('E41').Value = '  +0.38%  '
$ws.Range('E42').Value = '  -2.45%  '
$ws.Range('E43').Value = '  +1.77%  '
$ws.Range('E44').Value = '  +2.14%  '
$ws.Range('E45').Value = '  -4.26%  '
$ws.Range('E46').Value = '  +5.17%  '
$ws.Range('E47').Value = '  -0.59%  '
$ws.Range('E48').Value = '  +1.94%  '
$ws.Range('E49').Value = '  +0.07%  '
$ws.Range('E50').Value = '  +0.54%  '
$ws.Range('E51').Value = '  +0.52%  '

# Numeric-looking price updates: force text formatting so Excel does not
# coerce these into floating point numbers, then restore the original style
# so the cell keeps its original (unstyled) appearance.
$textCells = @{
    'D4' = '0.9986'
    'D5' = '236.52'
    'D6' = '0.9989'
    'D7' = '0.5177'
    'D8' = '40.38'
    'D9' = '0.2700'
    'D10' = '0.06211'
    'D12' = '0.06989'
    'D13' = '15.52'
    'D14' = '0.6396'
    'D15' = '4.489'
    'D16' = '78.12'
    'D17' = '0.9982'
    'D18' = '0.9987'
    'D20' = '11.67'
    'D21' = '0.000006704'
    'D23' = '4.076'
    'D24' = '8.363'
    'D25' = '5.190'
    'D26' = '136.57'
    'D27' = '1.482'
    'D28' = '1.836'
    'D29' = '15.17'
    'D30' = '103.23'
    'D31' = '0.08356'
    'D32' = '3.705'
    'D33' = '3.416'
    'D34' = '0.04402'
    'D35' = '2.644'
    'D36' = '0.9993'
    'D37' = '0.6079'
    'D38' = '2.726'
    'D39' = '0.01567'
    'D40' = '1.953'
    'D41' = '0.9990'
    'D42' = '102.46'
    'D43' = '0.3881'
    'D44' = '0.7492'
    'D45' = '4.931'
    'D46' = '0.05491'
    'D47' = '0.1110'
    'D48' = '6.079'
    'D49' = '30.28'
    'D50' = '52.76'
    'D51' = '1.002'
}
foreach ($cellRef in $textCells.Keys) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $textCells[$cellRef]
    $rng.Style = $origStyle
}

